# TrialsSetup 2026-02-11 12:00
# Update "Days remaining" figures on Sheet1 for two trials:
#   - REJOICE (MK-5909-003)  row 6: 9  -> 8
#   - REMASTER (CLOU)        row 8: 29 -> 28

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 8
$ws.Range("B8").Value = 28
